{"js": "// Update the date paragraph (first paragraph in the body) and the five\n// \"problem\" rows of the division-practice table. Each problem row is\n// every 4th row (0, 4, 8, 12, 16) of the single table in the document;\n// each has 5 cells. We replace each cell's text in place (using the\n// paragraph's Range so run/paragraph formatting such as font, size and\n// justification are preserved) rather than doing a blind global find\n// & replace, because several original cell values (e.g. \"59\u00f74=\",\n// \"18\u00f77=\") are duplicated across rows and must map to different new\n// values depending on position.\n\n// 1) Update the date line.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.getRange().insertText(\"2024-09-25 Wednesday\", \"Replace\");\n\n// 2) Update the practice-problem table cells.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index -> new cell values (left to right)\nconst rowUpdates = {\n  0: [\"24\u00f73=\", \"92\u00f72=\", \"23\u00f77=\", \"48\u00f77=\", \"94\u00f79=\"],\n  4: [\"30\u00f77=\", \"75\u00f72=\", \"39\u00f78=\", \"36\u00f79=\", \"78\u00f79=\"],\n  8: [\"72\u00f73=\", \"72\u00f79=\", \"24\u00f75=\", \"64\u00f75=\", \"21\u00f77=\"],\n  12: [\"52\u00f75=\", \"62\u00f76=\", \"15\u00f76=\", \"87\u00f74=\", \"51\u00f78=\"],\n  16: [\"79\u00f72=\", \"55\u00f73=\", \"73\u00f73=\", \"41\u00f74=\", \"93\u00f75=\"],\n};\n\nconst cellParagraphs = [];\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const values = rowUpdates[rowIndex];\n  for (let colIndex = 0; colIndex < values.length; colIndex++) {\n    const cell = table.getCell(Number(rowIndex), colIndex);\n    const cellParas = cell.body.paragraphs;\n    cellParas.load(\"items\");\n    cellParagraphs.push({ cellParas, text: values[colIndex] });\n  }\n}\nawait context.sync();\n\nfor (const { cellParas, text } of cellParagraphs) {\n  const para = cellParas.items[0];\n  para.getRange().insertText(text, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the date line and the five \"problem\" rows of the division\n# practice table. Each problem row is every 4th row (1, 5, 9, 13, 17 in\n# Word's 1-based indexing) of the single table in the document, each\n# with 5 cells. Cells are addressed by their (row, column) position --\n# not by searching for their old text -- because several original\n# values (e.g. \"59\u00f74=\", \"18\u00f77=\") are duplicated across rows and must\n# map to different new values depending on where they sit.\n#\n# Assigning plain text to Range.Text (rather than deleting/inserting)\n# keeps the existing run/paragraph formatting (font, size, alignment)\n# untouched, matching how Word itself behaves.\n\n$d = $word.ActiveDocument\n\n$d.Paragraphs(1).Range.Text = \"2024-09-25 Wednesday\"\n\n$t = $d.Tables(1)\n\n$t.Cell(1, 1).Range.Text = \"24\u00f73=\"\n$t.Cell(1, 2).Range.Text = \"92\u00f72=\"\n$t.Cell(1, 3).Range.Text = \"23\u00f77=\"\n$t.Cell(1, 4).Range.Text = \"48\u00f77=\"\n$t.Cell(1, 5).Range.Text = \"94\u00f79=\"\n\n$t.Cell(5, 1).Range.Text = \"30\u00f77=\"\n$t.Cell(5, 2).Range.Text = \"75\u00f72=\"\n$t.Cell(5, 3).Range.Text = \"39\u00f78=\"\n$t.Cell(5, 4).Range.Text = \"36\u00f79=\"\n$t.Cell(5, 5).Range.Text = \"78\u00f79=\"\n\n$t.Cell(9, 1).Range.Text = \"72\u00f73=\"\n$t.Cell(9, 2).Range.Text = \"72\u00f79=\"\n$t.Cell(9, 3).Range.Text = \"24\u00f75=\"\n$t.Cell(9, 4).Range.Text = \"64\u00f75=\"\n$t.Cell(9, 5).Range.Text = \"21\u00f77=\"\n\n$t.Cell(13, 1).Range.Text = \"52\u00f75=\"\n$t.Cell(13, 2).Range.Text = \"62\u00f76=\"\n$t.Cell(13, 3).Range.Text = \"15\u00f76=\"\n$t.Cell(13, 4).Range.Text = \"87\u00f74=\"\n$t.Cell(13, 5).Range.Text = \"51\u00f78=\"\n\n$t.Cell(17, 1).Range.Text = \"79\u00f72=\"\n$t.Cell(17, 2).Range.Text = \"55\u00f73=\"\n$t.Cell(17, 3).Range.Text = \"73\u00f73=\"\n$t.Cell(17, 4).Range.Text = \"41\u00f74=\"\n$t.Cell(17, 5).Range.Text = \"93\u00f75=\"\n"}
